# Update the "dSF" column (F) values on specific rows to reflect
# repulled data / recalculated means, per the commit:
#   "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -12
    4  = -3
    5  = 11
    9  = -3
    10 = -5
    12 = 4
    13 = -3
    21 = -1
    22 = -4
    24 = -4
    25 = -1
    34 = 3
    39 = -3
    42 = 3
    43 = 0
    52 = 0
    64 = -3
    69 = -8
    71 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
